$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: bump row height (matches the taller rendered row after formatting) ---
$ws.Rows(2).RowHeight = 18.75

# --- Refresh the Bearer Token for row 2 (new tokenId/iat; Lichess ID stays the same) ---
$ws.Range("E2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiZmY1N2FlYzAtYTkyZC00YWEzLTk5ZDYtYmExNWI4NzcyNDNiIiwiaWF0IjoxNzEzMzcwNzI2fQ.T2kMPItYfoLiJmZOKklaGt5h27O5Gszjnyeu2YgSOqY"

# --- Text columns (A, B, E, F, G) on row 2: explicit general alignment like the header row/columns ---
$ws.Range("A2").HorizontalAlignment = 1
$ws.Range("B2").HorizontalAlignment = 1
$ws.Range("E2").HorizontalAlignment = 1
$ws.Range("F2").HorizontalAlignment = 1
$ws.Range("G2").HorizontalAlignment = 1

# --- Numeric columns (C, D) on row 2: thousands-separated numbers, right aligned, automatic (theme) font color ---
$ws.Range("C1").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)
$ws.Range("C2:D2").Font.ThemeColor = 1
